$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the "weekly" record data which shifts down one row:
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), P (Precio $/Kg)
$cols = @(4, 10, 11, 12, 13, 16)

$firstRow = 135
$lastRow  = 180
$newRow   = 181
$lastCol  = 18

# 1) Snapshot the existing values for rows 135..180 before any writes,
#    so the shift-down can be applied without clobbering source data.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Build the new bottom row (181) as a cell-by-cell copy of row 180's
#    current content/format (this preserves number formats/styles
#    without introducing new style entries).
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item($lastRow, $c).Copy($ws.Cells.Item($newRow, $c))
}
$excel.CutCopyMode = $false

# 3) Shift the tracked columns down by one row: row r (136..181) takes
#    the snapshot that belonged to row r-1.
for ($r = $newRow; $r -ge ($firstRow + 1); $r--) {
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $snapshot["$($r - 1)-$c"]
    }
}

# 4) Row 135 becomes the new, latest weekly record.
$ws.Cells.Item($firstRow, 4).Value2 = 44468
$ws.Cells.Item($firstRow, 10).Value2 = 180

"done"
